$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Jianlei%Cao%NULL%0, Xiaorong%Hu%NULL%0, Wenlin%Cheng%NULL%0, Lei%Yu%NULL%0, Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0, Qiang%Liu%liuqiang@irm-cams.ac.cn%0, Qiang%Liu%liuqiang@irm-cams.ac.cn%0]'
$ws.Range("I2").Value = ''
$ws.Range("J2").Value = 'Springer Berlin Heidelberg'
$ws.Range("C3").Value = 'Unknown Title'
$ws.Range("D3").Value = 'Unknown Abstract'
$ws.Range("E3").Value = '[]'
$ws.Range("F3").Value = 'not found'
$ws.Range("G3").Value = 'N/A'
$ws.Range("H3").Value = '''1970-01-01'
$ws.Range("H3").Style = "Normal"
$ws.Range("E4").Value = '[Chaolin%Huang%NULL%0, Yeming%Wang%NULL%0, Xingwang%Li%NULL%0, Lili%Ren%NULL%0, Jianping%Zhao%NULL%0, Yi%Hu%NULL%0, Li%Zhang%NULL%0, Guohui%Fan%NULL%0, Jiuyang%Xu%NULL%0, Xiaoying%Gu%NULL%0, Zhenshun%Cheng%NULL%0, Ting%Yu%NULL%0, Jiaan%Xia%NULL%0, Yuan%Wei%NULL%0, Wenjuan%Wu%NULL%0, Xuelei%Xie%NULL%0, Wen%Yin%NULL%0, Hui%Li%NULL%0, Min%Liu%NULL%0, Yan%Xiao%NULL%0, Hong%Gao%NULL%0, Li%Guo%NULL%0, Jungang%Xie%NULL%0, Guangfa%Wang%NULL%0, Rongmeng%Jiang%NULL%0, Zhancheng%Gao%NULL%0, Qi%Jin%NULL%0, Jianwei%Wang%wangjw28@163.com%0, Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("I4").Value = ''
$ws.Range("J4").Value = 'Elsevier Ltd.'
$ws.Range("C5").Value = 'Unknown Title'
$ws.Range("D5").Value = 'Unknown Abstract'
$ws.Range("E5").Value = '[]'
$ws.Range("F5").Value = 'not found'
$ws.Range("G5").Value = 'N/A'
$ws.Range("H5").Value = '''1970-01-01'
$ws.Range("H5").Style = "Normal"
$ws.Range("E6").Value = '[Qiurong%Ruan%NULL%0, Kun%Yang%NULL%0, Kun%Yang%NULL%0, Wenxia%Wang%NULL%0, Wenxia%Wang%NULL%0, Lingyu%Jiang%NULL%0, Lingyu%Jiang%NULL%0, Jianxin%Song%songsingsjx@sina.com%0, Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("I6").Value = ''
$ws.Range("J6").Value = 'Springer Berlin Heidelberg'
$ws.Range("C7").Value = 'Unknown Title'
$ws.Range("D7").Value = 'Unknown Abstract'
$ws.Range("E7").Value = '[]'
$ws.Range("F7").Value = 'not found'
$ws.Range("G7").Value = 'N/A'
$ws.Range("I7").Value = ''
$ws.Range("C8").Value = 'Unknown Title'
$ws.Range("E8").Value = '[]'
$ws.Range("F8").Value = 'not found'
$ws.Range("G8").Value = 'N/A'
$ws.Range("I8").Value = ''
$ws.Range("C9").Value = 'Unknown Title'
$ws.Range("D9").Value = 'Unknown Abstract'
$ws.Range("E9").Value = '[]'
$ws.Range("F9").Value = 'not found'
$ws.Range("G9").Value = 'N/A'
$ws.Range("H9").Value = '''1970-01-01'
$ws.Range("H9").Style = "Normal"
$ws.Range("E10").Value = '[Xiaobo%Yang%NULL%0, Yuan%Yu%NULL%0, Jiqian%Xu%NULL%0, Huaqing%Shu%NULL%0, Jia''an%Xia%NULL%0, Hong%Liu%NULL%0, Yongran%Wu%NULL%0, Lu%Zhang%NULL%0, Zhui%Yu%NULL%0, Minghao%Fang%NULL%0, Ting%Yu%NULL%0, Yaxin%Wang%NULL%0, Shangwen%Pan%NULL%0, Xiaojing%Zou%NULL%0, Shiying%Yuan%NULL%0, You%Shang%NULL%0]'
$ws.Range("I10").Value = ''
$ws.Range("J10").Value = 'Elsevier Ltd.'
$ws.Range("E11").Value = '[Fei%Zhou%NULL%0, Ting%Yu%NULL%0, Ronghui%Du%NULL%0, Guohui%Fan%NULL%0, Ying%Liu%NULL%0, Zhibo%Liu%NULL%0, Jie%Xiang%NULL%0, Yeming%Wang%NULL%0, Bin%Song%NULL%0, Xiaoying%Gu%NULL%0, Lulu%Guan%NULL%0, Yuan%Wei%NULL%0, Hui%Li%NULL%0, Xudong%Wu%NULL%0, Jiuyang%Xu%NULL%0, Shengjin%Tu%NULL%0, Yi%Zhang%NULL%0, Hua%Chen%NULL%0, Bin%Cao%NULL%0]'
$ws.Range("I11").Value = ''
$ws.Range("J11").Value = 'Elsevier Ltd.'
